# Apply cryptos list update (auto-generated)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.038.95"
$ws.Range("E2").Value = "  +2.74%  "

$ws.Range("D3").Value = "1.652.19"
$ws.Range("E3").Value = "  +3.49%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'215.13"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("D8").Value = "'0.250"

$ws.Range("E9").Value = "  +1.62%  "

$ws.Range("D10").Value = "'19.84"
$ws.Range("E10").Value = "  +4.19%  "

$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").Value = "1.885.60"
$ws.Range("E12").Value = "  +3.55%  "

$ws.Range("D13").Value = "1.652.57"
$ws.Range("E13").Value = "  +3.35%  "

$ws.Range("D14").Value = "'4.08"
$ws.Range("E14").Value = "  +2.25%  "

$ws.Range("D15").Value = "'0.518"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").Value = "'65.32"
$ws.Range("E16").Value = "  +2.92%  "

$ws.Range("D17").Value = "'240.14"
$ws.Range("E17").Value = "  +4.27%  "

$ws.Range("D18").Value = "27.017.66"
$ws.Range("E18").Value = "  +2.71%  "

$ws.Range("E19").Value = "  +2.61%  "

$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  +4.49%  "

$ws.Range("E23").Value = "  +3.34%  "

$ws.Range("E24").Value = "  +3.45%  "

$ws.Range("D25").Value = "'145.65"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("E28").Value = "  +1.85%  "

$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("E32").Value = "  +3.40%  "

$ws.Range("D33").Value = "1.519.84"
$ws.Range("E33").Value = "  +1.04%  "

$ws.Range("D34").Value = "'3.07"
$ws.Range("E34").Value = "  +5.13%  "

$ws.Range("E35").Value = "  +8.64%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").Value = "'0.580"
$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'0.886"
$ws.Range("E38").Value = "  +8.52%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0169"
$ws.Range("E39").Value = "  +2.95%  "

$ws.Range("E40").Value = "  +2.74%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  +4.37%  "

$ws.Range("D43").Value = "'65.93"
$ws.Range("E43").Value = "  +9.00%  "

$ws.Range("D44").Value = "1.792.14"

$ws.Range("D45").Value = "'0.774"
$ws.Range("E45").Value = "  +2.12%  "

$ws.Range("E46").Value = "  -2.93%  "

$ws.Range("D47").Value = "'89.53"
$ws.Range("E47").Value = "  +1.24%  "

$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("E51").Value = "  +2.19%  "

